# This script reproduces the commit that inserts two new price rows
# (row 304 and row 305) into the "Limón" sheet, pushing every existing
# row from 304-395 down by two positions (to 306-397). The workbook's
# dimension grows from A1:T395 to A1:T397.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 304. Doing this twice at the same
# index shifts the original row 304 (and everything below it) down by
# two rows in total, landing the old row 304 at row 306.
$ws.Rows.Item(304).Insert()
$ws.Rows.Item(304).Insert()

# ---- New row 304 ----
$ws.Range("A304").Value = 7
$ws.Range("B304").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C304").Value = "Ñuble"
$ws.Range("D304").Value = 44463
$ws.Range("E304").Value = 16
$ws.Range("F304").Value = "Fruta"
$ws.Range("G304").Value = 100102
$ws.Range("H304").Value = "Cítricos"
$ws.Range("I304").Value = 100102003
$ws.Range("J304").Value = "Limón"
$ws.Range("K304").Value = "Sin especificar"
$ws.Range("L304").Value = "1a amarillo"
$ws.Range("M304").Value = 240
$ws.Range("N304").Value = 4800
$ws.Range("O304").Value = 5000
$ws.Range("P304").Value = 4900
$ws.Range("Q304").Value = "$/malla 18 kilos"
$ws.Range("R304").Value = "Región Metropolitana"
$ws.Range("S304").Value = 272
$ws.Range("T304").Value = 18

# ---- New row 305 ----
$ws.Range("A305").Value = 7
$ws.Range("B305").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C305").Value = "Ñuble"
$ws.Range("D305").Value = 44463
$ws.Range("E305").Value = 16
$ws.Range("F305").Value = "Fruta"
$ws.Range("G305").Value = 100102
$ws.Range("H305").Value = "Cítricos"
$ws.Range("I305").Value = 100102003
$ws.Range("J305").Value = "Limón"
$ws.Range("K305").Value = "Sin especificar"
$ws.Range("L305").Value = "2a amarillo"
$ws.Range("M305").Value = 240
$ws.Range("N305").Value = 3800
$ws.Range("O305").Value = 4200
$ws.Range("P305").Value = 4000
$ws.Range("Q305").Value = "$/malla 18 kilos"
$ws.Range("R305").Value = "Región Metropolitana"
$ws.Range("S305").Value = 222
$ws.Range("T305").Value = 18
